$wb = $excel.ActiveWorkbook

# ---- Overview ----
$ws = $wb.Worksheets.Item('Overview')
$ws.Range('A1').Value = 'File Name'
$ws.Range('B1').Value = 'zh-cn'
$ws.Range('C1').Value = 'de-de'
$ws.Range('D1').Value = 'Latest Handoff Date'
$ws.Range('B2').Value = 'Handed back: in sync with en-US'
$ws.Range('C2').Value = 'Handed back: in sync with en-US'
$ws.Range('D2').Value = '2016-03-22 21:12:45'
$ws.Range('B3').Value = 'Handed back: in sync with en-US'
$ws.Range('C3').Value = 'Handed back: in sync with en-US'
$ws.Range('D3').Value = '2016-03-22 21:11:02'
$ws.Range('B4').Value = 'Handed back: in sync with en-US'
$ws.Range('C4').Value = 'Handed back: in sync with en-US'
$ws.Range('D4').Value = '2016-03-22 21:11:02'
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md', '', '', 'faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md', '', '', 'ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/ffffffab0b8d44-08cb-4112-bb04-99628528bfad.md', '', '', 'ffffffab0b8d44-08cb-4112-bb04-99628528bfad.md')

# ---- zh-cn ----
$ws = $wb.Worksheets.Item('zh-cn')
$ws.Range('A1').Value = 'Source File Name'
$ws.Range('B1').Value = 'File Extension'
$ws.Range('C1').Value = 'Status'
$ws.Range('D1').Value = 'Latest Handoff File'
$ws.Range('E1').Value = 'Latest Handoff Datetime'
$ws.Range('F1').Value = 'Latest Target File'
$ws.Range('G1').Value = 'Latest Handback File'
$ws.Range('H1').Value = 'Latest Handback DateTime'
$ws.Range('I1').Value = 'Reference Tokens'
$ws.Range('J1').Value = 'Handoff Reason'
$ws.Range('K1').Value = 'Dependency From'
$ws.Range('L1').Value = 'Error Detail'
$ws.Range('B2').Value = '.md'
$ws.Range('C2').Value = 'Handed back: in sync with en-US'
$ws.Range('E2').Value = '2016-03-22 21:12:41'
$ws.Range('H2').Value = '2016-03-22 21:13:04'
$ws.Range('J2').Value = 'Include'
$ws.Range('B3').Value = '.md'
$ws.Range('C3').Value = 'Handed back: in sync with en-US'
$ws.Range('E3').Value = '2016-03-22 21:10:58'
$ws.Range('H3').Value = '2016-03-22 21:11:29'
$ws.Range('J3').Value = 'Include'
$ws.Range('B4').Value = '.md'
$ws.Range('C4').Value = 'Handed back: in sync with en-US'
$ws.Range('E4').Value = '2016-03-22 21:10:58'
$ws.Range('H4').Value = '2016-03-22 21:11:29'
$ws.Range('J4').Value = 'Include'
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md', '', '', 'faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md')
$ws.Hyperlinks.Add($ws.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ffc014fffac0fcdc0be71b90c0c84f2e7abe4d68/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.zh-cn.xlf', '', '', 'faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md', '', '', 'faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md')
$ws.Hyperlinks.Add($ws.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ffc014fffac0fcdc0be71b90c0c84f2e7abe4d68/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.zh-cn.xlf', '', '', 'faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md', '', '', 'ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md')
$ws.Hyperlinks.Add($ws.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3087839134feca2713bf27c7a424e7afc32d48ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2d082ab0d7e9b915974610d76271fcab73e77c4b/e2e/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md')
$ws.Hyperlinks.Add($ws.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/036358e3a4fbec3e1b786289aa15d712cb8a9d45/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/ffffffab0b8d44-08cb-4112-bb04-99628528bfad.md', '', '', 'ffffffab0b8d44-08cb-4112-bb04-99628528bfad.md')
$ws.Hyperlinks.Add($ws.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3087839134feca2713bf27c7a424e7afc32d48ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2d082ab0d7e9b915974610d76271fcab73e77c4b/e2e/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md')
$ws.Hyperlinks.Add($ws.Range('G4'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/036358e3a4fbec3e1b786289aa15d712cb8a9d45/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf')

# ---- de-de ----
$ws = $wb.Worksheets.Item('de-de')
$ws.Range('A1').Value = 'Source File Name'
$ws.Range('B1').Value = 'File Extension'
$ws.Range('C1').Value = 'Status'
$ws.Range('D1').Value = 'Latest Handoff File'
$ws.Range('E1').Value = 'Latest Handoff Datetime'
$ws.Range('F1').Value = 'Latest Target File'
$ws.Range('G1').Value = 'Latest Handback File'
$ws.Range('H1').Value = 'Latest Handback DateTime'
$ws.Range('I1').Value = 'Reference Tokens'
$ws.Range('J1').Value = 'Handoff Reason'
$ws.Range('K1').Value = 'Dependency From'
$ws.Range('L1').Value = 'Error Detail'
$ws.Range('B2').Value = '.md'
$ws.Range('C2').Value = 'Handed back: in sync with en-US'
$ws.Range('E2').Value = '2016-03-22 21:12:45'
$ws.Range('H2').Value = '2016-03-22 21:13:10'
$ws.Range('J2').Value = 'Include'
$ws.Range('B3').Value = '.md'
$ws.Range('C3').Value = 'Handed back: in sync with en-US'
$ws.Range('E3').Value = '2016-03-22 21:11:02'
$ws.Range('H3').Value = '2016-03-22 21:11:38'
$ws.Range('J3').Value = 'Include'
$ws.Range('B4').Value = '.md'
$ws.Range('C4').Value = 'Handed back: in sync with en-US'
$ws.Range('E4').Value = '2016-03-22 21:11:02'
$ws.Range('H4').Value = '2016-03-22 21:11:38'
$ws.Range('J4').Value = 'Include'
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md', '', '', 'faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md')
$ws.Hyperlinks.Add($ws.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68af4d6a6e5c63750b8ba4b5bd9137ccf67e553e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.de-de.xlf', '', '', 'faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md', '', '', 'faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md')
$ws.Hyperlinks.Add($ws.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68af4d6a6e5c63750b8ba4b5bd9137ccf67e553e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.de-de.xlf', '', '', 'faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md', '', '', 'ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md')
$ws.Hyperlinks.Add($ws.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f69f789c581b2b476ae44bb0ae79d3d67f47e62/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8d7885a453560134b037738f0376603a0c245c6d/e2e/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md')
$ws.Hyperlinks.Add($ws.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/c26ef30b4b06c1237078a348596e274fd545df6a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/1304274c9f4c20980a404a180c60b0d3c3e59aec/e2e/ffffffab0b8d44-08cb-4112-bb04-99628528bfad.md', '', '', 'ffffffab0b8d44-08cb-4112-bb04-99628528bfad.md')
$ws.Hyperlinks.Add($ws.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f69f789c581b2b476ae44bb0ae79d3d67f47e62/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8d7885a453560134b037738f0376603a0c245c6d/e2e/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md')
$ws.Hyperlinks.Add($ws.Range('G4'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/c26ef30b4b06c1237078a348596e274fd545df6a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf', '', '', 'c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf')
